$d = $word.ActiveDocument
$d.Content.Find.Execute("ATRIBUIÇÕES", $true, $false, $false, $false, $false,
                         $true, 1, $false, "PROCESSOS", 2)
